# Apply updates to the "Metadata" sheet (sheet1) and leave the
# "Include from Health Data Conn" sheet (sheet2) values unchanged
# (only their underlying shared-string indices shift because new
# strings were inserted earlier in the table - Excel handles that
# bookkeeping automatically when we just set .Value on cells).

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 used to be "Contact" / "No display for ContactDetail" (duplicated
# on row 11 too). It becomes "Jurisdiction" / "United States of America",
# and the duplicate row is removed entirely (rows shift up by one).
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Delete what was row 11 (the duplicate Contact row), shifting everything
# below it up by one row.
$meta.Rows.Item(11).Delete()
